# Commit: "Added config2 for Hawkeye"
#
# Config2 (Hawkeye config) previously borrowed its Warmup/Sim "50 mil"
# placeholder columns and cross-referenced Config1's row 5 data for the
# Hawkeye policy. This change gives Config2 its own measured row 5
# (Hawkeye) and row 6 (OPTGen) data -- mirroring the layout already used
# in Config1 -- and drops the now-unused "Warmup"/"50 mil"/"Sim" helper
# cells from both sheets' header/data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Config1: drop the leftover "Warmup"/"Sim" header labels (K2/L2)
# and the "50 mil" placeholder cells in rows 3-4 (K3:L3, K4:L4).
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config1")

$ws1.Range("K2").ClearContents()
$ws1.Range("L2").ClearContents()

$ws1.Range("K3").ClearContents()
$ws1.Range("L3").ClearContents()

$ws1.Range("K4").ClearContents()
$ws1.Range("L4").ClearContents()

# ---------------------------------------------------------------
# Config2: same header cleanup -- J2 keeps its text but becomes the
# "Hit Rate" label (matching Config1's J2 after the shared-string
# table is compacted); K2 is cleared and L2 becomes a new, empty,
# identically-styled cell.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config2")

$null = $ws2.Range("K2").Copy($ws2.Range("L2"))
$ws2.Range("K2").ClearContents()
$ws2.Range("L2").ClearContents()
$ws2.Range("J2").Value = "Hit Rate"

# drop "50 mil" placeholders in rows 3-4 (J3:K3, J4:K4)
$ws2.Range("J3").ClearContents()
$ws2.Range("K3").ClearContents()

$ws2.Range("J4").ClearContents()
$ws2.Range("K4").ClearContents()

# Row 5 (Hawkeye): give Config2 its own measured data instead of
# referencing back into Config1.
$ws2.Range("C5").Value = 50000002
$ws2.Range("D5").Value = 131405787
$ws2.Range("E5").Value = 2517423
$ws2.Range("F5").Value = 2374604
$ws2.Range("G5").Value = 142819
$ws2.Range("H5").Formula = "=(C5/D5)"
$ws2.Range("I5").Formula = "=G5/(C5/1000)"

# Row 6 (OPTGen): newly populated with measured data (was previously
# blank and resolved to #DIV/0!).
$ws2.Range("C6").Value = 50000002
$ws2.Range("D6").Value = 131405787
$ws2.Range("E6").Value = 34143
$ws2.Range("F6").Value = 43676
$ws2.Range("G6").Formula = "=E6-F6"
$ws2.Range("H6").Formula = "=(C6/D6)"
$ws2.Range("I6").Formula = "=G6/(C6/1000)"
$ws2.Range("J6").Formula = "=F6/E6"

# Selection ends on Config2!C7, but Config1 stays the active
# (tab-selected) sheet, matching its pre-existing C7 selection.
$null = $ws2.Range("C7").Select()
$null = $ws1.Activate()
